# Actualización automática 2025-12-01 08:30:07
#
# "VENTAS POR GRUPO": the PORCELANATO (M) / PIEDRA SINTERIZADA (L) sales
# that had fallen into this month's columns are cleared back to 0, and the
# "x de 18" counters on row 20 for those two columns follow suit.
#
# "VENTA MENSUAL": the rolling 4-month window shifts forward by one month
# (agosto/septiembre/octubre/noviembre -> septiembre/octubre/noviembre/
# diciembre), so every data column C..F shifts left by one (new C = old D,
# new D = old E, new E = old F, new F = 0 for the freshly-opened month).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("M2").Value = 0
$ws1.Range("M9").Value = 0
$ws1.Range("M13").Value = 0
$ws1.Range("L18").Value = 0
$ws1.Range("L20").Value = "0 de 18"
$ws1.Range("M20").Value = "0 de 18"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# column widths follow the same left-shift as the data (old D/E/F widths
# move into C/D/E; F keeps its width since it's now the "empty" month).
# Excel's ColumnWidth property (character units) round-trips through the
# default-font padding (+5/6 of a character) before it lands in the raw
# OOXML <col width="..."> attribute, so back the padding out here to land
# on the intended raw widths of 16 / 13 / 15.
$ws2.Columns.Item(3).ColumnWidth = 16 - 0.8333333333333334
$ws2.Columns.Item(4).ColumnWidth = 13 - 0.8333333333333334
$ws2.Columns.Item(5).ColumnWidth = 15 - 0.8333333333333334

# month headers
$ws2.Range("C1").Value = "septiembre"
$ws2.Range("D1").Value = "octubre"
$ws2.Range("E1").Value = "noviembre"
$ws2.Range("F1").Value = "diciembre"

# data rows: shift C<-D, D<-E, E<-F, F<-0 for every row with figures
$rows = 2..20
foreach ($r in $rows) {
    $oldD = $ws2.Cells.Item($r, 4).Value()
    $oldE = $ws2.Cells.Item($r, 5).Value()
    $oldF = $ws2.Cells.Item($r, 6).Value()

    $ws2.Cells.Item($r, 3).Value = $oldD
    $ws2.Cells.Item($r, 4).Value = $oldE
    $ws2.Cells.Item($r, 5).Value = $oldF
    $ws2.Cells.Item($r, 6).Value = 0
}
